# Applies the "Now also does reversed cars" edit:
#  - Renames the "inter (wrong)" legend label (shared string) to "count"
#    and repurposes it as the header for a new column L that counts
#    wheel/intra/inter-axle detections per sensor group, with a SUM total.
#  - Removes the old stray J10/K10 legend row ("inter (wrong)" / 738).
#  - Updates the recorded timestamps for the tail of train 1 (rows 53:56)
#    and replaces the whole of train 2's timestamps (rows 58:108) with the
#    newly-recorded "reversed cars" pass data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column L ("count") ------------------------------------------------
# Renames shared string index 7 from "inter (wrong)" -> "count" (it is
# currently only used by K10) and uses it as the L1 header.
$ws.Range("K10").Value = "count"
$ws.Range("L1").Value = "count"

$ws.Range("L2").Value = 24
$ws.Range("L3").Value = 4
$ws.Range("L5").Value = 2
$ws.Range("L6").Value = 12
$ws.Range("L8").Value = 2
$ws.Range("L9").Value = 11
$ws.Range("L11").Formula = "=SUM(L2:L9)"

# Drop the old stray "inter (wrong)" legend entry on row 10.
$ws.Range("J10").ClearContents()
$ws.Range("K10").ClearContents()

# --- Timestamp corrections for the tail of train 1 (rows 53:56) ------------
$ws.Range("A53").Value = 50771
$ws.Range("A54").Value = 51131
$ws.Range("A55").Value = 52146
$ws.Range("A56").Value = 52506

# --- Train 2 ("reversed cars") timestamps, rows 58:108 ---------------------
$ws.Range("A58").Value = 11812
$ws.Range("A59").Value = 12187
$ws.Range("A60").Value = 12783
$ws.Range("A61").Value = 12826
$ws.Range("A62").Value = 15048
$ws.Range("A63").Value = 15329
$ws.Range("A64").Value = 15916
$ws.Range("A65").Value = 18222
$ws.Range("A66").Value = 18525
$ws.Range("A67").Value = 19094
$ws.Range("A68").Value = 19125
$ws.Range("A69").Value = 21352
$ws.Range("A70").Value = 21647
$ws.Range("A71").Value = 22229
$ws.Range("A72").Value = 22518
$ws.Range("A73").Value = 24558
$ws.Range("A74").Value = 24902
$ws.Range("A75").Value = 25436
$ws.Range("A76").Value = 25546
$ws.Range("A77").Value = 27708
$ws.Range("A78").Value = 28614
$ws.Range("A79").Value = 28759
$ws.Range("A80").Value = 30889
$ws.Range("A81").Value = 31222
$ws.Range("A82").Value = 31754
$ws.Range("A83").Value = 32092
$ws.Range("A84").Value = 34052
$ws.Range("A85").Value = 34363
$ws.Range("A86").Value = 34934
$ws.Range("A87").Value = 35227
$ws.Range("A88").Value = 37247
$ws.Range("A89").Value = 37512
$ws.Range("A90").Value = 38097
$ws.Range("A91").Value = 38495
$ws.Range("A92").Value = 40388
$ws.Range("A93").Value = 40695
$ws.Range("A94").Value = 41267
$ws.Range("A95").Value = 41597
$ws.Range("A96").Value = 43567
$ws.Range("A97").Value = 43825
$ws.Range("A98").Value = 44435
$ws.Range("A99").Value = 44780
$ws.Range("A100").Value = 46683
$ws.Range("A101").Value = 47005
$ws.Range("A102").Value = 47615
$ws.Range("A103").Value = 47839
$ws.Range("A104").Value = 49889
$ws.Range("A105").Value = 50250
$ws.Range("A106").Value = 50794
$ws.Range("A107").Value = 51162
$ws.Range("A108").Value = 52167

# --- Window/selection state --------------------------------------------
$ws.Range("M13").Select()
